$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 0.721220700610068
$ws.Range("E4").Value = 0.703415187125182
$ws.Range("F4").Value = 0.759706650286704
$ws.Range("G4").Value = 0.575101001313061
$ws.Range("H4").Value = 0.731116681894394
$ws.Range("I4").Value = 0.683377424047149
$ws.Range("J4").Value = 0.659704988708263
$ws.Range("K4").Value = 0.759539408850822
$ws.Range("L4").Value = 0.690381293922693
$ws.Range("M4").Value = 0.721853509181885
$ws.Range("N4").Value = 0.618731104894089
$ws.Range("B5").Value = 0.716369766863604
$ws.Range("C5").Value = 0.772717497952248
$ws.Range("D5").Value = 0.768126880551711
$ws.Range("E5").Value = 0.754584861941685
$ws.Range("F5").Value = 0.87256696213353
$ws.Range("G5").Value = 0.844461583387313
$ws.Range("H5").Value = 0.841389015556442
$ws.Range("I5").Value = 0.65867391008924
$ws.Range("J5").Value = 0.658259163542218
$ws.Range("K5").Value = 0.702279283302204
$ws.Range("L5").Value = 0.770751009217637
$ws.Range("M5").Value = 0.925808443317861
$ws.Range("N5").Value = 0.561959539989434
$ws.Range("D6").Value = 0.703425024778734
$ws.Range("E6").Value = 0.690418272415403
$ws.Range("F6").Value = 0.817088824971887
$ws.Range("G6").Value = 0.707832398863728
$ws.Range("H6").Value = 0.743080145118413
$ws.Range("I6").Value = 0.682786565268408
$ws.Range("J6").Value = 0.629978960602202
$ws.Range("K6").Value = 0.559364498889858
$ws.Range("M6").Value = 0.729122438104651
$ws.Range("N6").Value = 0.667115559981111
$ws.Range("D7").Value = 0.716516257894719
$ws.Range("E7").Value = 0.707070890580701
$ws.Range("F7").Value = 0.816121829740342
$ws.Range("G7").Value = 0.639318978823409
$ws.Range("H7").Value = 0.771177034328568
$ws.Range("I7").Value = 0.693178534315937
$ws.Range("J7").Value = 0.571035672573194
$ws.Range("K7").Value = 0.556955266701761
$ws.Range("N7").Value = 0.669890735490999
$ws.Range("D8").Value = 0.438889632437241
$ws.Range("E8").Value = 0.441403577339026
$ws.Range("F8").Value = 0.692723367340905
$ws.Range("H8").Value = 0.511112243128561
$ws.Range("I8").Value = 0.456721110351558
$ws.Range("N8").Value = 0.402399638605333
$ws.Range("D9").Value = 0.433525987240537
$ws.Range("E9").Value = 0.622061067971373
$ws.Range("F9").Value = 0.764034268173467
$ws.Range("G9").Value = 0.627378598315257
$ws.Range("H9").Value = 0.710473338089272
$ws.Range("I9").Value = 0.580152462807448
$ws.Range("J9").Value = 0.529014769509267
$ws.Range("K9").Value = 0.585754346694377
$ws.Range("L9").Value = 0.568451357363088
$ws.Range("M9").Value = 0.880923373592755
$ws.Range("N9").Value = 0.553829259151759
